$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Add 14 new "T1 - NIST sphere N" header columns (AD1:AQ1), copying
#    the same formatting (bold / border / centered) used by the other
#    header cells such as Z1.
# ------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("AD1:AQ1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headers = @(
    "T1 - NIST sphere 1",
    "T1 - NIST sphere 2",
    "T1 - NIST sphere 3",
    "T1 - NIST sphere 4",
    "T1 - NIST sphere 5",
    "T1 - NIST sphere 6",
    "T1 - NIST sphere 7",
    "T1 - NIST sphere 8",
    "T1 - NIST sphere 9",
    "T1 - NIST sphere 10",
    "T1 - NIST sphere 11",
    "T1 - NIST sphere 12",
    "T1 - NIST sphere 13",
    "T1 - NIST sphere 14"
)

# Column 30 = AD, ... Column 43 = AQ
$startCol = 30
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# ------------------------------------------------------------------
# 2) Fix the casing of the "t1map" -> "T1map" token inside every
#    NIFTI filename in column D (rows 2-57).
# ------------------------------------------------------------------
for ($row = 2; $row -le 57; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $old = $cell.Value2
    if ($old -ne $null -and $old.Contains("t1map.nii.gz")) {
        $new = $old.Replace("t1map.nii.gz", "T1map.nii.gz")
        $cell.Value = $new
    }
}
